# Refresh the cached "datetimeFigureOut" auto-date field shown on the
# Slide Master and every Slide Layout's Date placeholder (PowerPoint
# recomputes/re-caches this text whenever the deck is opened/saved on a
# later date; here we simulate that resave by updating the cached text
# from 3/11/2017 to 3/5/2018).

$p = $ppt.ActivePresentation
$newDate = "3/5/2018"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Type -eq 14) {
            $ph = $sh.PlaceholderFormat
            if ($ph.Type -eq 16) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# Slide Master
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every Custom Layout hanging off the master
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}
